$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Actions Required" column now holds the google-form voting note ---
$ws.Range("E4").Value = "Carlos: Add Ideas to repository, create Google Form for voting; Everyone: Vote on google form(s) until we have a final 5"

# --- Row 5: new meeting entry for 1/26 ---
$ws.Range("B5").Value = "Select Project from list of 5, as narrowed down from the google form surveys"

$ws.Range("C5").Value = "Carlos, Courtnie, Lucas, Patrick"

$ws.Range("D5").Value = "Chose to turn in Fire Alarm and Teaching Assistant for Divide and Conquer - Initial Document due 2/3"
$ws.Range("D5").WrapText = $true

$ws.Range("F5").Value = "We decided to Use Fire Alarm and Teaching Assistant; See Deadlines spreadsheet for details;"
$ws.Range("F5").WrapText = $true

$e5Text = "Please complete the following by 2/1: Carlos: Project Milestones; Courtnie: House of Quality, Engineering Specifications; Lucas: budget and financing, goals and objectives;  Patrick: Project Block Diagram, Decision Matrix;"
$ws.Range("E5").Value = $e5Text
$ws.Range("E5").WrapText = $true
$e5 = $ws.Range("E5")
# "Please complete the following by 2/1:" -> italic
$e5.Characters(1, 37).Font.Italic = $true
# "Carlos" -> bold
$e5.Characters(39, 6).Font.Bold = $true
# "Courtnie" -> bold
$e5.Characters(67, 8).Font.Bold = $true
# "Lucas:" -> bold
$e5.Characters(123, 6).Font.Bold = $true
# "Patrick:" -> bold
$e5.Characters(175, 8).Font.Bold = $true

# Row 5 grew taller to fit the new wrapped content
$ws.Rows("5:5").RowHeight = 60

# --- Update the "last updated" banner (B1) ---
$ws.Range("B1").Value = "Last updated: 01/26/2017 8:35 PM by Carlos"

# Selection moved to B1 (matches the saved sheetView)
$ws.Range("B1").Select()
